$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 49.74535073696344, -70.77256502277163),
    @(3, 83.27769771251813, -49.91575236160403),
    @(4, 103.5734216329357, 18.76143815029367),
    @(5, 115.1890028215457, 33.22021346205078),
    @(6, 118.5723899384148, 35.47684965769099),
    @(7, 121.9043418767182, 38.31892265686311),
    @(8, 126.4788201136039, 44.19965757108333),
    @(9, 145.0084388579184, 60.32743001907394),
    @(10, 169.6849164300639, 72.1203844943892),
    @(11, 204.1283344512795, 76.86995262199872),
    @(12, 60.93090702498469, -101.2949112474405),
    @(13, 228.5251664856048, -124.1229380778302),
    @(14, 276.0308352822454, -113.3030278155987),
    @(15, 270.2913324481655, -43.70583623240577),
    @(16, 259.6606186505904, -25.48782393053351),
    @(17, 252.0147344005555, -18.85575500459782),
    @(18, 255.8251879225983, -18.2470456333632),
    @(19, 165.9977784716485, 5.957775829676551),
    @(20, 142.966670473727, -4.258449207339063),
    @(21, 145.1914362708839, 13.48412920588432),
    @(22, -97.74890414797861, 6.55103626039878),
    @(23, -116.8385183846027, 44.72900005582461),
    @(24, -132.1571509619447, 87.53925227821099),
    @(25, -152.6872065836371, 96.79474626138233),
    @(26, -138.7850698755932, 140.5605615001747),
    @(27, -133.21259527928, 155.8618242334191),
    @(28, -133.9124958948067, 155.3899253307704),
    @(29, -134.1000625002073, 155.6522068960797),
    @(30, -135.1738009553617, 157.6066559315938),
    @(31, -132.3585532804421, 165.9649555558985),
    @(32, -236.9009260781619, -138.458473601824),
    @(33, -203.1988613098107, -108.9593938900891),
    @(34, -189.6615886520328, -80.54484495394192),
    @(35, -188.5293796166974, -79.64785988736745),
    @(36, -195.0929860585928, -98.45056163988015),
    @(37, -192.8312958570538, -95.9266941830513),
    @(38, -189.2144447399565, -102.8962588314725),
    @(39, -181.7298584576943, -99.56758521885411),
    @(40, -180.7318146759262, -99.68374674921844),
    @(41, -230.1318686931855, -91.28739448357405)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
